$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 15: "др. культуры" -> "др. культуры:" and clear unit cells B15/C15
$ws.Range("A15").Value = "др. культуры:"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""

# Row 16: add new "овсюг" entry in column A, clear unit cells F16/G16
$ws.Range("A16").Value = "овсюг"
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = ""

# Update the active selection to A16
[void]$ws.Range("A16").Select()
